$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "Back and Forw Jump"
$ws.Range("C17").Value = "Xavier Trillo"

$ws.Range("B16").Value = "Neutral Jump Smooth"
$ws.Range("C16").Value = "Xavier Trillo"
$ws.Range("D16").Value = "0 hours and 30 minutes"
$ws.Range("E16").Value = "0 hours and 45 minutes"

$ws.Range("D17").Value = "1 hour "
$ws.Range("E17").Value = "1 hour and 45 minutes"

$ws.Range("D18").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
